# Auto-generated Excel COM-interop script to apply cryptos.xlsx price/volume update
# (Wed Mar 15 19:47:43 UTC 2023 GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5 and 6: BNB / USDC swapped position (coin name, link, price, volume)
$ws.Range("B5").Value = "USDC"
$ws.Range("C5").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"

# Rows 40 and 41: FraxShare / TrustWalletToken swapped position (coin name, link, price, volume)
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"

# Update Price (column D) and Volume(1h) (column E) for rows 2-51.
# Price values are stored as plain text (not numbers) in the sheet; values that
# look numeric are written with a leading "'" quote-prefix so Excel keeps them
# as text, then the style is reset to Normal so no extra text-format style is
# introduced (matches the original inline-string cells, which carry no style).
$ws.Range("D2").Value = "24.344.59"
$ws.Range("E2").Value = "  -1.92%  "
$ws.Range("D3").Value = "1.646.78"
$ws.Range("E3").Value = "  -4.05%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.59%  "
$ws.Range("D5").Value = "'1.001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").Value = "'306.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.87%  "
$ws.Range("D7").Value = "'0.3619"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.22%  "
$ws.Range("D8").Value = "'47.33"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.62%  "
$ws.Range("D9").Value = "'0.3269"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.32%  "
$ws.Range("D10").Value = "'1.110"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.90%  "
$ws.Range("D11").Value = "'0.06885"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.85%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").Value = "'5.916"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.76%  "
$ws.Range("D14").Value = "'19.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.79%  "
$ws.Range("D15").Value = "1.645.81"
$ws.Range("E15").Value = "  -4.44%  "
$ws.Range("D16").Value = "'6.522"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.58%  "
$ws.Range("D17").Value = "'0.00001039"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.80%  "
$ws.Range("D18").Value = "'0.06482"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.94%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("D20").Value = "'76.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.23%  "
$ws.Range("D21").Value = "'5.877"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.84%  "
$ws.Range("D22").Value = "'15.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -9.90%  "
$ws.Range("D23").Value = "'12.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -10.32%  "
$ws.Range("D24").Value = "24.346.67"
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("D25").Value = "'2.430"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("D26").Value = "'2.283"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -18.58%  "
$ws.Range("D27").Value = "'145.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.54%  "
$ws.Range("D28").Value = "'18.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -10.27%  "
$ws.Range("D29").Value = "1.836.98"
$ws.Range("E29").Value = "  -3.95%  "
$ws.Range("D30").Value = "'123.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.95%  "
$ws.Range("D31").Value = "'1.147"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.22%  "
$ws.Range("D32").Value = "'4.041"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.61%  "
$ws.Range("D33").Value = "'5.517"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -19.21%  "
$ws.Range("D34").Value = "'0.08309"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.45%  "
$ws.Range("D35").Value = "'1.673"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.96%  "
$ws.Range("D36").Value = "'12.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -11.58%  "
$ws.Range("D37").Value = "'5.162"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.47%  "
$ws.Range("D38").Value = "'0.05999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.24%  "
$ws.Range("D39").Value = "'0.02198"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.82%  "
$ws.Range("D40").Value = "'1.201"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.75%  "
$ws.Range("D41").Value = "'8.171"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.15%  "
$ws.Range("D42").Value = "'0.2029"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.98%  "
$ws.Range("D43").Value = "'0.9997"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").Value = "'0.5785"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -10.29%  "
$ws.Range("D45").Value = "'3.712"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.57%  "
$ws.Range("D46").Value = "'12.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.91%  "
$ws.Range("D47").Value = "'0.5513"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.22%  "
$ws.Range("D48").Value = "'120.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.57%  "
$ws.Range("D49").Value = "'1.921"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -10.57%  "
$ws.Range("D50").Value = "'0.06871"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.58%  "
$ws.Range("D51").Value = "'73.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.77%  "
